$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the instruction counts (B6:B10) and data cache counts (B14:B15)
$ws.Range("B6").Value = 3383
$ws.Range("B7").Value = 218
$ws.Range("B8").Value = 833
$ws.Range("B9").Value = 500
$ws.Range("B10").Value = 512
$ws.Range("B14").Value = 300
$ws.Range("B15").Value = 292

# Clear the leftover "Grades" helper table (headers stay, but are now blank; data rows removed)
$ws.Range("K20").Value = $null
$ws.Range("L20").Value = $null
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null

$ws.Range("K21:K24").EntireRow.Delete()

# Move selection to reflect the new active cell
$ws.Range("B15").Select()
